$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New rolling 89-day window: oldest day (2025-10-10) dropped, newest day
# (2026-01-07) appended. Column A dates shift forward by one day; Column C
# item counts are the refreshed GSC export figures for the new window.
$dates = @("2025-10-11","2025-10-12","2025-10-13","2025-10-14","2025-10-15","2025-10-16","2025-10-17","2025-10-18","2025-10-19","2025-10-20","2025-10-21","2025-10-22","2025-10-23","2025-10-24","2025-10-25","2025-10-26","2025-10-27","2025-10-28","2025-10-29","2025-10-30","2025-10-31","2025-11-01","2025-11-02","2025-11-03","2025-11-04","2025-11-05","2025-11-06","2025-11-07","2025-11-08","2025-11-09","2025-11-10","2025-11-11","2025-11-12","2025-11-13","2025-11-14","2025-11-15","2025-11-16","2025-11-17","2025-11-18","2025-11-19","2025-11-20","2025-11-21","2025-11-22","2025-11-23","2025-11-24","2025-11-25","2025-11-26","2025-11-27","2025-11-28","2025-11-29","2025-11-30","2025-12-01","2025-12-02","2025-12-03","2025-12-04","2025-12-05","2025-12-06","2025-12-07","2025-12-08","2025-12-09","2025-12-10","2025-12-11","2025-12-12","2025-12-13","2025-12-14","2025-12-15","2025-12-16","2025-12-17","2025-12-18","2025-12-19","2025-12-20","2025-12-21","2025-12-22","2025-12-23","2025-12-24","2025-12-25","2025-12-26","2025-12-27","2025-12-28","2025-12-29","2025-12-30","2025-12-31","2026-01-01","2026-01-02","2026-01-03","2026-01-04","2026-01-05","2026-01-06","2026-01-07")
$items = @(13,26,34,41,49,50,59,63,66,72,81,81,83,84,85,90,83,90,93,92,101,112,115,108,105,101,95,87,82,76,54,47,43,41,38,35,31,29,26,26,25,26,26,25,25,27,27,27,27,27,27,27,27,26,25,25,25,26,26,27,29,29,30,30,31,31,32,31,31,32,32,32,32,30,31,32,32,28,28,28,28,30,29,28,28,27,27,27,27)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $dateCell = $ws.Cells.Item($row, 1)
    # Force text so the date-like string isn't reinterpreted as a serial date,
    # then drop the number-format override so the cell keeps the sheet's
    # default (unstyled) appearance.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dates[$i]
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 3).Value = $items[$i]
}
